$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a typo in the "pagi" (morning) time-slot label ---
# Row 13 held "06:55 - 8:00" which should read "06:55 - 7:00".
$ws.Range("A13").Value = "06:55 – 7:00"

# --- Update existing pagi (morning) counts for rows 2-7 ---
$b257 = @(17, 27, 33, 27, 46, 52)
$c257 = @(2, 0, 3, 1, 2, 0)
for ($i = 0; $i -lt $b257.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $b257[$i]
    $ws.Cells.Item($row, 3).Value = $c257[$i]
}

# --- Add newly-recorded siang (midday, rows 8-25) counts ---
$bNew = @(43, 57, 134, 90, 38, 130, 95, 132, 100, 123, 80, 72, 111, 143, 102, 95, 169, 154)
$cNew = @(2, 1, 4, 4, 3, 5, 3, 4, 2, 5, 2, 1, 4, 2, 5, 3, 6, 3)
for ($i = 0; $i -lt $bNew.Length; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 2).Value = $bNew[$i]
    $ws.Cells.Item($row, 3).Value = $cNew[$i]
}

# Center-align the newly-populated data cells (rows 8-25), matching the
# look of the already-filled rows above them.
$ws.Range("B8:C25").HorizontalAlignment = -4108 # xlCenter

# --- Touch (but leave blank) the still-empty rows 26-49 so they pick up ---
# --- the same centred number formatting as the rest of the table.      ---
$ws.Range("B26:C49").HorizontalAlignment = -4108 # xlCenter

# --- Re-center the TOTAL row's computed sums ---
$ws.Range("B50:C50").HorizontalAlignment = -4108 # xlCenter

# Recalculate so B50/C50 (SUM formulas) reflect the new figures.
$excel.Calculate()

# --- Restore the saved view/selection state ---
$excel.ActiveWindow.ScrollRow = 20
$ws.Range("B52").Select()
